$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45181 -> 45182) for every data row, rows 2 through 78.
for ($r = 2; $r -le 78; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
